# Daily attendance processing - 2026-01-27 22:00:31
# Swap the order of "Recorded By" names in column G from "System, <email>"
# to "<email>, System" for every row in the used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System, ")) {
        $rest = $val.Substring(8)
        $cell.Value2 = "$rest, System"
    }
}
